# Insert a new weekly price record as row 97 on the "Orégano" sheet.
# Every existing row from 97 downward (97-180) shifts down by one
# (97->98, ..., 180->181), which Rows(97).Insert() handles natively,
# matching Excel's own "insert row, push down" semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(97).Insert()

$row = 97
$ws.Cells.Item($row, 1).Value  = 6
$ws.Cells.Item($row, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($row, 3).Value  = "Metropolitana"
$ws.Cells.Item($row, 4).Value  = 44658
$ws.Cells.Item($row, 5).Value  = 13
$ws.Cells.Item($row, 6).Value  = 100112029
$ws.Cells.Item($row, 7).Value  = "Orégano"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 35
$ws.Cells.Item($row, 11).Value = 13000
$ws.Cells.Item($row, 12).Value = 14000
$ws.Cells.Item($row, 13).Value = 13457
$ws.Cells.Item($row, 14).Value = "`$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 4486
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
